$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.334.05"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.90%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.566.84"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.09%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.000"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "289.16"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.73%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3756"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.77%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.30"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.71%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3420"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.168"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.72%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07654"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.01%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.39"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.019"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.69%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.946"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.48%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001135"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.24%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.556.71"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.58%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.07"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.58%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06713"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.68%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.244"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.92%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.59"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.98%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5263"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.17%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.95"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.25%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "22.328.92"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.98%  "

$ws.Range("E26").Value = "  -1.52%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.797"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.52%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.14"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.54%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "145.63"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.970"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "125.50"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.83%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.731.75"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.59%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.026"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.73%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.235"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -9.94%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.006"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.38%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.10"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -9.77%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08498"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02534"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.13%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2323"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.20%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.544"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.78%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.316"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.55%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.06401"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.33%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.74"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.79%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6391"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.61%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.18"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.34%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9997"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5991"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.33%  "

$ws.Range("E48").Value = "  -4.39%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.097"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.91%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.272"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.37%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "124.35"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.13%  "
